# Update TPM-derived values in the Adam9-Itgb5 LR-pairs worksheet
# following re-computation with new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.627806666666667
$ws.Range("H2").Value = 10.88342
$ws.Range("I2").Value = 0.08036675778141429
$ws.Range("J2").Value = 0.08036675778141429
$ws.Range("M2").Value = 8.033114333333334
$ws.Range("N2").Value = 24.099343
$ws.Range("O2").Value = 0.1374088679258946
$ws.Range("P2").Value = 0.1374088679258946
$ws.Range("Q2").Value = 29.14258573256222
$ws.Range("R2").Value = 262.28327159306
$ws.Range("S2").Value = 0.01104310520561872
$ws.Range("T2").Value = 0.01104310520561872
$ws.Range("G3").Value = 3.627806666666667
$ws.Range("H3").Value = 10.88342
$ws.Range("I3").Value = 0.08036675778141429
$ws.Range("J3").Value = 0.08036675778141429
$ws.Range("O3").Value = 0.6355200716780686
$ws.Range("P3").Value = 0.6355200716780686
$ws.Range("Q3").Value = 134.7853195590733
$ws.Range("R3").Value = 1213.06787603166
$ws.Range("S3").Value = 0.05107468766577838
$ws.Range("T3").Value = 0.05107468766577838
$ws.Range("G4").Value = 3.627806666666667
$ws.Range("H4").Value = 10.88342
$ws.Range("I4").Value = 0.08036675778141429
$ws.Range("J4").Value = 0.08036675778141429
$ws.Range("M4").Value = 13.27489133333333
$ws.Range("N4").Value = 39.824674
$ws.Range("O4").Value = 0.2270710603960369
$ws.Range("P4").Value = 0.2270710603960369
$ws.Range("Q4").Value = 48.15873927834223
$ws.Range("R4").Value = 433.42865350508
$ws.Range("S4").Value = 0.01824896491001719
$ws.Range("T4").Value = 0.01824896491001719
$ws.Range("I5").Value = 0.6506403335968259
$ws.Range("J5").Value = 0.6506403335968259
$ws.Range("M5").Value = 8.033114333333334
$ws.Range("N5").Value = 24.099343
$ws.Range("O5").Value = 0.1374088679258946
$ws.Range("P5").Value = 0.1374088679258946
$ws.Range("Q5").Value = 235.9351332111771
$ws.Range("R5").Value = 2123.416198900594
$ws.Range("S5").Value = 0.08940375166646625
$ws.Range("T5").Value = 0.08940375166646626
$ws.Range("I6").Value = 0.6506403335968259
$ws.Range("J6").Value = 0.6506403335968259
$ws.Range("O6").Value = 0.6355200716780686
$ws.Range("P6").Value = 0.6355200716780686
$ws.Range("S6").Value = 0.4134949914440972
$ws.Range("T6").Value = 0.4134949914440972
$ws.Range("I7").Value = 0.6506403335968259
$ws.Range("J7").Value = 0.6506403335968259
$ws.Range("M7").Value = 13.27489133333333
$ws.Range("N7").Value = 39.824674
$ws.Range("O7").Value = 0.2270710603960369
$ws.Range("P7").Value = 0.2270710603960369
$ws.Range("Q7").Value = 389.8877975752991
$ws.Range("R7").Value = 3508.990178177692
$ws.Range("S7").Value = 0.1477415904862624
$ws.Range("T7").Value = 0.1477415904862624
$ws.Range("G8").Value = 12.14251133333333
$ws.Range("H8").Value = 36.427534
$ws.Range("I8").Value = 0.2689929086217598
$ws.Range("J8").Value = 0.2689929086217598
$ws.Range("M8").Value = 8.033114333333334
$ws.Range("N8").Value = 24.099343
$ws.Range("O8").Value = 0.1374088679258946
$ws.Range("P8").Value = 0.1374088679258946
$ws.Range("Q8").Value = 97.54218183446245
$ws.Range("R8").Value = 877.8796365101621
$ws.Range("S8").Value = 0.03696201105380963
$ws.Range("T8").Value = 0.03696201105380963
$ws.Range("G9").Value = 12.14251133333333
$ws.Range("H9").Value = 36.427534
$ws.Range("I9").Value = 0.2689929086217598
$ws.Range("J9").Value = 0.2689929086217598
$ws.Range("O9").Value = 0.6355200716780686
$ws.Range("P9").Value = 0.6355200716780686
$ws.Range("Q9").Value = 451.1354712892647
$ws.Range("R9").Value = 4060.219241603382
$ws.Range("S9").Value = 0.170950392568193
$ws.Range("T9").Value = 0.170950392568193
$ws.Range("G10").Value = 12.14251133333333
$ws.Range("H10").Value = 36.427534
$ws.Range("I10").Value = 0.2689929086217598
$ws.Range("J10").Value = 0.2689929086217598
$ws.Range("M10").Value = 13.27489133333333
$ws.Range("N10").Value = 39.824674
$ws.Range("O10").Value = 0.2270710603960369
$ws.Range("P10").Value = 0.2270710603960369
$ws.Range("Q10").Value = 161.1905184637685
$ws.Range("R10").Value = 1450.714666173916
$ws.Range("S10").Value = 0.06108050499975725
$ws.Range("T10").Value = 0.06108050499975725
